$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Hoja 1 - Compresion Table")
$ws2.Range("F4").Value = 73.22
$ws2.Range("F4").NumberFormat = "General"
$ws2.Range("B5").Value = 79.34999999999999
$ws2.Range("B5").NumberFormat = "General"
$ws2.Range("C5").Value = 76.98999999999999
$ws2.Range("C5").NumberFormat = "General"
$ws2.Range("D5").Value = 76.76000000000001
$ws2.Range("D5").NumberFormat = "General"
$ws2.Range("E5").Value = 76.86
$ws2.Range("E5").NumberFormat = "General"
$ws2.Range("F5").Value = 76.91
$ws2.Range("F5").NumberFormat = "General"
$ws2.Range("B6").Value = 68.48999999999999
$ws2.Range("B6").NumberFormat = "General"
$ws2.Range("C6").Value = 67.69
$ws2.Range("C6").NumberFormat = "General"
$ws2.Range("D6").Value = 66.59
$ws2.Range("D6").NumberFormat = "General"
$ws2.Range("E6").Value = 66.18000000000001
$ws2.Range("E6").NumberFormat = "General"
$ws2.Range("F6").Value = 66.18000000000001
$ws2.Range("F6").NumberFormat = "General"

$ws3 = $wb.Worksheets.Item("Hoja 1 - Compresion Time in Sec")
$ws3.Range("F4").Value = 211.76
$ws3.Range("F4").NumberFormat = "General"
$ws3.Range("B5").Value = 212.83
$ws3.Range("B5").NumberFormat = "General"
$ws3.Range("C5").Value = 212.31
$ws3.Range("C5").NumberFormat = "General"
$ws3.Range("D5").Value = 212.1
$ws3.Range("D5").NumberFormat = "General"
$ws3.Range("E5").Value = 212.23
$ws3.Range("E5").NumberFormat = "General"
$ws3.Range("F5").Value = 213.5
$ws3.Range("F5").NumberFormat = "General"
$ws3.Range("B6").Value = 215.4
$ws3.Range("B6").NumberFormat = "General"
$ws3.Range("C6").Value = 216.46
$ws3.Range("C6").NumberFormat = "General"
$ws3.Range("D6").Value = 215.04
$ws3.Range("D6").NumberFormat = "General"
$ws3.Range("E6").Value = 215.01
$ws3.Range("E6").NumberFormat = "General"
$ws3.Range("F6").Value = 215.63
$ws3.Range("F6").NumberFormat = "General"

$ws4 = $wb.Worksheets.Item("Hoja 1 - Decompresion Time in S")
$ws4.Range("F4").Value = 29.2
$ws4.Range("F4").NumberFormat = "General"
$ws4.Range("B5").Value = 27.09
$ws4.Range("B5").NumberFormat = "General"
$ws4.Range("C5").Value = 27.89
$ws4.Range("C5").NumberFormat = "General"
$ws4.Range("D5").Value = 28.23
$ws4.Range("D5").NumberFormat = "General"
$ws4.Range("E5").Value = 28.15
$ws4.Range("E5").NumberFormat = "General"
$ws4.Range("F5").Value = 28.41
$ws4.Range("F5").NumberFormat = "General"
$ws4.Range("B6").Value = 28.1
$ws4.Range("B6").NumberFormat = "General"
$ws4.Range("C6").Value = 29.48
$ws4.Range("C6").NumberFormat = "General"
$ws4.Range("D6").Value = 29.4
$ws4.Range("D6").NumberFormat = "General"
$ws4.Range("E6").Value = 29.8
$ws4.Range("E6").NumberFormat = "General"
$ws4.Range("F6").Value = 29.83
$ws4.Range("F6").NumberFormat = "General"
